$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D cells to Text format so numeric-looking price
# strings (e.g. "0.9997", "28.199.27") are preserved exactly as text and not
# auto-converted/rounded into floating point numbers by Excel.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '28.199.27'
$ws.Range('E2').Value = '  +3.81%  '
$ws.Range('D3').Value = '1.810.35'
$ws.Range('E3').Value = '  +1.62%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.76%  '
$ws.Range('D5').Value = '338.81'
$ws.Range('E5').Value = '  +0.86%  '
$ws.Range('D6').Value = '0.9997'
$ws.Range('E6').Value = '  -0.51%  '
$ws.Range('D7').Value = '0.3931'
$ws.Range('E7').Value = '  +3.87%  '
$ws.Range('D8').Value = '0.3497'
$ws.Range('E8').Value = '  +1.99%  '
$ws.Range('D9').Value = '48.00'
$ws.Range('E9').Value = '  -1.09%  '
$ws.Range('D10').Value = '1.176'
$ws.Range('E10').Value = '  -1.28%  '
$ws.Range('D11').Value = '0.07558'
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('D12').Value = '0.9997'
$ws.Range('E12').Value = '  -0.89%  '
$ws.Range('D13').Value = '22.12'
$ws.Range('E13').Value = '  +1.40%  '
$ws.Range('D14').Value = '6.517'
$ws.Range('E14').Value = '  +1.33%  '
$ws.Range('D15').Value = '1.809.42'
$ws.Range('E15').Value = '  +1.15%  '
$ws.Range('D16').Value = '7.151'
$ws.Range('E16').Value = '  +1.37%  '
$ws.Range('D17').Value = '0.00001104'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').Value = '0.06725'
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('D19').Value = '85.41'
$ws.Range('E19').Value = '  +1.14%  '
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').Value = '17.75'
$ws.Range('E21').Value = '  +2.72%  '
$ws.Range('D22').Value = '6.574'
$ws.Range('E22').Value = '  +0.78%  '
$ws.Range('D23').Value = '28.191.81'
$ws.Range('E23').Value = '  +3.52%  '
$ws.Range('D24').Value = '12.40'
$ws.Range('E24').Value = '  -0.35%  '
$ws.Range('D25').Value = '2.394'
$ws.Range('E25').Value = '  -1.61%  '
$ws.Range('D26').Value = '21.52'
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range('D27').Value = '1.483'
$ws.Range('E27').Value = '  -0.94%  '
$ws.Range('D28').Value = '2.528'
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('D29').Value = '154.29'
$ws.Range('E29').Value = '  +0.93%  '
$ws.Range('D30').Value = '2.012.05'
$ws.Range('E30').Value = '  +0.99%  '
$ws.Range('D31').Value = '136.43'
$ws.Range('E31').Value = '  +2.25%  '
$ws.Range('D32').Value = '6.237'
$ws.Range('E32').Value = '  +3.47%  '
$ws.Range('D33').Value = '4.010'
$ws.Range('E33').Value = '  -1.11%  '
$ws.Range('D34').Value = '0.08854'
$ws.Range('E34').Value = '  +2.04%  '
$ws.Range('D35').Value = '13.32'
$ws.Range('E35').Value = '  +1.87%  '
$ws.Range('D36').Value = '0.02442'
$ws.Range('E36').Value = '  +4.81%  '
$ws.Range('D37').Value = '0.6933'
$ws.Range('E37').Value = '  +1.21%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '5.456'
$ws.Range('E38').Value = '  +0.31%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.06522'
$ws.Range('E39').Value = '  +2.33%  '
$ws.Range('B40').Value = 'WEMIXTOKEN'
$ws.Range('C40').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D40').Value = '1.607'
$ws.Range('E40').Value = '  -2.60%  '
$ws.Range('D41').Value = '0.2215'
$ws.Range('E41').Value = '  +1.26%  '
$ws.Range('D42').Value = '1.260'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').Value = '8.551'
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('D44').Value = '14.56'
$ws.Range('E44').Value = '  +0.26%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.6437'
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('B46').Value = 'Frax'
$ws.Range('C46').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D46').Value = '1.000'
$ws.Range('E46').Value = '  -0.26%  '
$ws.Range('D47').Value = '3.876'
$ws.Range('E47').Value = '  +0.64%  '
$ws.Range('D48').Value = '2.151'
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('D49').Value = '131.18'
$ws.Range('E49').Value = '  +1.60%  '
$ws.Range('D50').Value = '0.07183'
$ws.Range('E50').Value = '  -0.18%  '
$ws.Range('D51').Value = '80.29'
$ws.Range('E51').Value = '  +1.28%  '

# Restore original (default/general) formatting now that the text values are set.
$ws.Range('D2:D51').ClearFormats()
